$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9-23 down to 10-24
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with data (same as surrounding rows, but with new values)
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44781
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 100112017
$ws.Cells.Item(9, 7).Value = "Ramas de apio"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 40
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 5000
$ws.Cells.Item(9, 14).Value = "$/paquete"
$ws.Cells.Item(9, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(9, 16).Value = 5000
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Copy the date style from row 10 (D column) to the new row 9 D cell
$ws.Cells.Item(10, 4).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4122)
